$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F4").Value = 343
$ws.Range("F6").Value = 730
$ws.Range("F9").Value = 2403
$ws.Range("F15").Value = 5873
$ws.Range("F16").Value = 91
$ws.Range("F17").Value = 1139
$ws.Range("F19").Value = 1333
$ws.Range("F23").Value = 1723
$ws.Range("F24").Value = 254
$ws.Range("F26").Value = 613
$ws.Range("F29").Value = 4159
$ws.Range("F33").Value = 3587
$ws = $wb.Worksheets.Item(2)
$ws.Range("F8").Value = 10
$ws.Range("F21").Value = 567
$ws.Range("F22").Value = 207
$ws.Range("F23").Value = 341
$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 384
$ws.Range("F13").Value = 409
$ws.Range("F14").Value = 1065
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 384
$ws.Range("F6").Value = 343
$ws.Range("F11").Value = 730
$ws.Range("F15").Value = 2403
$ws.Range("F23").Value = 5873
$ws.Range("F24").Value = 409
$ws.Range("F25").Value = 1139
$ws.Range("F28").Value = 1723
$ws.Range("F29").Value = 254
$ws.Range("F32").Value = 613
$ws.Range("F34").Value = 4161
$ws.Range("F37").Value = 3587
